$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Kayitlar")
$ws5 = $wb.Worksheets.Item("Merkez İlçe")

$ws1.Rows.Item(1107).Delete()
$ws5.Rows.Item(568).Delete()
